$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 560
$ws.Range("I13").Value = 425
$ws.Range("J13").Value = 650
$ws.Range("K13").Value = 425
$ws.Range("L13").Value = 650
$ws.Range("M13").Value = -256
$ws.Range("N13").Value = -988
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 5000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 5000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -5936
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 5000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -5468
$ws.Range("H34").Value = 3914.375
$ws.Range("I34").Value = 188
$ws.Range("J34").Value = 29999
$ws.Range("K34").Value = 188
$ws.Range("L34").Value = 29999
$ws.Range("M34").Value = 15
$ws.Range("N34").Value = -30405
$ws.Range("H36").Value = 3914.375
$ws.Range("I36").Value = 188
$ws.Range("J36").Value = 29999
$ws.Range("K36").Value = 188
$ws.Range("L36").Value = 29999
$ws.Range("M36").Value = 527
$ws.Range("N36").Value = -31429
$ws.Range("H40").Value = 2833.6667
$ws.Range("I40").Value = 2001
$ws.Range("J40").Value = 3250
$ws.Range("K40").Value = 2001
$ws.Range("L40").Value = 3250
$ws.Range("M40").Value = -1826
$ws.Range("N40").Value = -3600
$ws.Range("H43").Value = 2250.5
$ws.Range("J43").Value = 2250.5
$ws.Range("L43").Value = 2250.5
$ws.Range("N43").Value = -2388.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 944.4286
$ws.Range("I16").Value = 935.1667
$ws.Range("K16").Value = 935.1667
$ws.Range("M16").Value = -648.1667
$ws.Range("H60").Value = 20359
$ws.Range("J60").Value = 27931.666
$ws.Range("L60").Value = 27931.666
$ws.Range("N60").Value = -28953.666
$ws.Range("H99").Value = 5395.35
$ws.Range("I99").Value = 3638.8572
$ws.Range("J99").Value = 6341.154
$ws.Range("K99").Value = 3638.8572
$ws.Range("L99").Value = 6341.154
$ws.Range("M99").Value = -2140.8572
$ws.Range("N99").Value = -9337.154
$ws.Range("H113").Value = 944.4286
$ws.Range("I113").Value = 935.1667
$ws.Range("K113").Value = 935.1667
$ws.Range("M113").Value = 1234.8333
$ws.Range("H126").Value = 5395.35
$ws.Range("I126").Value = 3638.8572
$ws.Range("J126").Value = 6341.154
$ws.Range("K126").Value = 10916.5716
$ws.Range("L126").Value = 19023.462
$ws.Range("M126").Value = -8446.571599999999
$ws.Range("N126").Value = -23963.462
$ws.Range("H134").Value = 1912
$ws.Range("I134").Value = 1744.1428
$ws.Range("K134").Value = 5232.428400000001
$ws.Range("M134").Value = -2697.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9470
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H80").Value = 8213.429
$ws.Range("J80").Value = 14998
$ws.Range("L80").Value = 44994
$ws.Range("N80").Value = -46866
$ws.Range("H83").Value = 8213.429
$ws.Range("J83").Value = 14998
$ws.Range("L83").Value = 134982
$ws.Range("N83").Value = -144342
$ws.Range("H132").Value = 1259.8
$ws.Range("I132").Value = 1462.5
$ws.Range("J132").Value = 1124.6666
$ws.Range("K132").Value = 13162.5
$ws.Range("L132").Value = 10121.9994
$ws.Range("M132").Value = -10632.5
$ws.Range("N132").Value = -15181.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 387.4762
$ws.Range("I97").Value = 381.4375
$ws.Range("J97").Value = 406.8
$ws.Range("K97").Value = 381.4375
$ws.Range("L97").Value = 406.8
$ws.Range("M97").Value = 114.5625
$ws.Range("N97").Value = -1398.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2532.25
$ws.Range("I22").Value = 1439.5
$ws.Range("K22").Value = 1439.5
$ws.Range("M22").Value = -1144.5
$ws.Range("H27").Value = 2532.25
$ws.Range("I27").Value = 1439.5
$ws.Range("K27").Value = 1439.5
$ws.Range("M27").Value = -1332.5
$ws.Range("H31").Value = 4608.222
$ws.Range("I31").Value = 2684.375
$ws.Range("J31").Value = 19999
$ws.Range("K31").Value = 2684.375
$ws.Range("L31").Value = 19999
$ws.Range("M31").Value = -2436.375
$ws.Range("N31").Value = -20495
$ws.Range("H46").Value = 1908.0312
$ws.Range("J46").Value = 2739.4666
$ws.Range("L46").Value = 2739.4666
$ws.Range("N46").Value = -3115.4666
$ws.Range("H55").Value = 734.375
$ws.Range("I55").Value = 155
$ws.Range("J55").Value = 817.1429000000001
$ws.Range("K55").Value = 155
$ws.Range("L55").Value = 817.1429000000001
$ws.Range("M55").Value = 18
$ws.Range("N55").Value = -1163.1429
$ws.Range("H82").Value = 3999.889
$ws.Range("J82").Value = 3999.889
$ws.Range("L82").Value = 3999.889
$ws.Range("N82").Value = -4721.889
$ws.Range("H85").Value = 3999.889
$ws.Range("J85").Value = 3999.889
$ws.Range("L85").Value = 3999.889
$ws.Range("N85").Value = -6495.889
$ws.Range("H94").Value = 15000
$ws.Range("J94").Value = 15000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352
$ws.Range("H132").Value = 13071
$ws.Range("I132").Value = 13427.714
$ws.Range("J132").Value = 12714.286
$ws.Range("K132").Value = 40283.142
$ws.Range("L132").Value = 38142.858
$ws.Range("M132").Value = -37753.142
$ws.Range("N132").Value = -43202.858
$ws.Range("H136").Value = 4234.4
$ws.Range("I136").Value = 3529
$ws.Range("J136").Value = 4939.8
$ws.Range("K136").Value = 10587
$ws.Range("L136").Value = 14819.4
$ws.Range("M136").Value = -8037
$ws.Range("N136").Value = -19919.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3558.1667
$ws.Range("I132").Value = 3269.8
$ws.Range("K132").Value = 9809.400000000001
$ws.Range("M132").Value = -7279.400000000001
